$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'306.11"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'6.54%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'32.09"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'8.38%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.273"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'2.73%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.07356"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'9.70%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'7.841"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'6.88%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'3.791"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'11.50%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'1.488"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'8.42%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.9119"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'-0.93%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.01653"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'2,461.54%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.1697"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'6.49%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.07554"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'11.60%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.08036"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'3.98%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.03027"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'3.37%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.09998"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'11.21%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.001506"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-5.07%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.04567"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'1.88%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'0.006192"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'-1.42%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'3.477"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'0.94%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'2.233"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'0.29%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.3284"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'2.14%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.1346"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'3.33%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'4.336"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'6.76%"
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'4.23%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.001232"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'3.22%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.004441"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'7.73%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'0.0001314"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'9.72%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'0.0001757"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'8.73%"
$ws.Range("E28").Style = "Normal"
$ws.Range("D40").Value = "'0.04505"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'5.44%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.007317"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'9.07%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1352"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'9.01%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.002285"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'5.01%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.01407"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'17.47%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.00006084"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'6.67%"
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'-4.13%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.01312"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'0.51%"
$ws.Range("E47").Style = "Normal"
